$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain text, matching the source data
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.276.75'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '3.500.30'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '586.80'
$ws.Range("E5").Value = '  +1.37%  '
$ws.Range("D6").Value = '134.15'
$ws.Range("E6").Value = '  +3.19%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("E9").Value = '  +2.34%  '
$ws.Range("D10").Value = '7.22'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("D11").Value = '0.384'
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '4.092.66'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("E14").Value = '  +3.09%  '
$ws.Range("D15").Value = '3.493.21'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '64.308.88'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '25.69'
$ws.Range("E17").Value = '  -5.14%  '
$ws.Range("D18").Value = '9.97'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '5.79'
$ws.Range("E19").Value = '  +3.21%  '
$ws.Range("D20").Value = '13.78'
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("D21").Value = '387.18'
$ws.Range("E21").Value = '  +1.31%  '
$ws.Range("D22").Value = '0.568'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '3.638.78'
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").Value = '74.16'
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("D27").Value = '0.0000114'
$ws.Range("E27").Value = '  +3.38%  '
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '8.30'
$ws.Range("E30").Value = '  +2.12%  '
$ws.Range("D31").Value = '2.24'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("E32").Value = '  -5.23%  '
$ws.Range("D33").Value = '3.520.04'
$ws.Range("E33").Value = '  +0.92%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.150'
$ws.Range("E34").Value = '  +4.69%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '23.58'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").Value = '5.21'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("D39").Value = '6.90'
$ws.Range("E39").Value = '  +1.13%  '
$ws.Range("D40").Value = '163.29'
$ws.Range("E40").Value = '  -2.43%  '
$ws.Range("D41").Value = '0.0782'
$ws.Range("E41").Value = '  -1.22%  '
$ws.Range("D42").Value = '0.806'
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").Value = '25.98'
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = '41.78'
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("D46").Value = '4.42'
$ws.Range("E46").Value = '  +1.83%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '1.18'
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '1.65'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Value = '2.482.31'
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").Value = '0.900'
$ws.Range("E51").Value = '  +2.30%  '
